$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = 5
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 22

# Row 6
$ws.Range("F6").Value = 100
$ws.Range("G6").Value = 21

# Row 7
$ws.Range("D7").Value = 9
$ws.Range("F7").Value = 100
$ws.Range("G7").Value = 20

# Row 8
$ws.Range("F8").Value = 100

# Row 9
$ws.Range("F9").Value = 100
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = 22

# Row 13
$ws.Range("D13").Value = 10
$ws.Range("F13").Value = 84
$ws.Range("G13").Value = 15

# Row 14
$ws.Range("D14").Value = 10
$ws.Range("F14").Value = 84
$ws.Range("G14").Value = 15

# Row 15
$ws.Range("D15").Value = 10
$ws.Range("F15").Value = 84
$ws.Range("G15").Value = 15

# Row 16
$ws.Range("D16").Value = 10
$ws.Range("F16").Value = 84
$ws.Range("G16").Value = 15

# Row 17
$ws.Range("D17").Value = 10.2
$ws.Range("F17").Value = 91
$ws.Range("G17").Value = 13

# Row 19
$ws.Range("D19").Value = 10.2
$ws.Range("F19").Value = 91
$ws.Range("G19").Value = 13

# Row 20
$ws.Range("D20").Value = 10
$ws.Range("F20").Value = 91
$ws.Range("G20").Value = 13

# Row 21
$ws.Range("D21").Value = 10
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 13

# Row 22
$ws.Range("D22").Value = 10
$ws.Range("F22").Value = 91
$ws.Range("G22").Value = 13

# Row 23
$ws.Range("D23").Value = 10
$ws.Range("F23").Value = 91
$ws.Range("G23").Value = 13

# Row 24
$ws.Range("D24").Value = 10
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 13

# Row 25
$ws.Range("D25").Value = 10
$ws.Range("F25").Value = 91
$ws.Range("G25").Value = 13

# Update the selected cell to match final workbook state
$ws.Range("J23").Select()
